$d = $word.ActiveDocument

# New paragraph texts (final document has 37 paragraphs; originally 49).
$newTexts = @(
    '**Original Document:**',
    'THIS AGREEMENT dated ${agreement-date}, 2025, by and between Omdena Inc, a Delaware Corporation (Hereinafter called "Party1") and ______________________________(Hereinafter called "Party2"), registered address at ______________________________________________, with registration number ${register0number}. (only if a company).',
    'WHEREAS, Party 1 and Party 2, for their mutual benefit and pursuant to a working relationship which has been or may be established, anticipate that either party may disclose or deliver to the working relationship which has been or may be established, documents, components, parts, information, drawings, data, sketches, plans, programs, specifications, techniques, processes, software, inventions, and other materials, both written and oral, of a secret, confidential or proprietary nature, including without limitation any information relating to marketing, finance, forecasts, invention, research, design, or development of information system and any supportive or incidental subsystems, and any subject matter claimed in or disclosed by any patent application prepared or filed by or behalf of either party, in any jurisdiction, and any amendments or supplements thereto (collectively, "Proprietary Information"); and',
    'WHEREAS, both Party 1 and Party 2 desire to assure that the confidentiality of any Proprietary Information is maintained;',
    'NOW, THEREFORE, in consideration of the foregoing premises and the mutual covenants contained herein, Party 1 and Party 2 hereby agree as follows:',
    '1. For a period of thirty-six (36) months from the date hereof, either party shall hold in trust and confidence, and not disclose to others or use for own benefit or for the benefit of another, any Proprietary Information that is disclosed to either parties at any time between the date hereof and twelve (12) months thereafter. Both parties can disclose Proprietary Information received under this Agreement to person within its organization only if such persons (i) have a need to know and (ii) are bound in writing to protect the confidentiality of such Proprietary Information. This paragraph 1 shall survive and continue after any expiration or termination of this Agreement shall bind both parties, its employees, agents, representatives, successors, heirs and assigns.',
    '2. The undertakings and obligations of either party under this Agreement shall not apply to any Proprietary Information which: (a) is described in an issued patent anywhere in the world, is disclosed in a printed publication available to the public, or is otherwise in the public domain through no action or fault of Party 2; (b) is generally disclosed to third parties by either party without restriction on such third parties, or is approved for release by written authorization of Party 1; (c) if not designated "confidential" at the time of first disclosure hereunder, or is not later designated in writing by either party within thirty (30) days from disclosure to be of a secret, confidential or proprietary nature; or (d) is shown to the receiving party within ten (10) days from disclosure, by underlying documentation to have been known by the receiving party before receipt of the confidential information and/or to have been developed by the receiving party completely independent of any disclosure.',
    '3. Title to all property received by the receiving party, including all Proprietary Information, shall remain at all times the sole property of the disclosing party, and this Agreement shall not be construed to grant to the receiving party any patents, licenses, or similar rights to such property and Proprietary Information disclosed to the receiving party hereunder.',
    '4. Receiving party shall, upon request of the disclosing party, return all documents, drawings and other tangible materials, including all Proprietary Information and all manifestations thereof, delivered by the disclosing party, and all copies and reproductions thereof.',
    '5. Non-Circumvention and Non-Solicitation Either party agrees not to circumvent the services and relationships. For six months from signing this agreement, either party must not submit proposals or solicit offers directly from contacts introduced by the disclosing party.',
    '6. The parties further agree to the following terms and conditions:',
    'i. Any breach by either party''s obligations under this Agreement will result in an irreparable inquiry for which damages and other legal remedies will be inadequate. In seeking enforcement of any of these obligations, the disclosing party will be entitled (in addition to other remedies) to preliminary and permanent injunctive and other equitable relief to prevent, discontinue and/or restrain the breach of this Agreement.',
    'ii. If any provision of this Agreement is invalid or unenforceable, then such provision shall be construed and limited to the extent necessary, or severed if necessary, in order to eliminate such invalidity or unenforceability, and the other provisions of this Agreement shall not be affected thereby.',
    'iii. In any dispute over whether information or matter is Proprietary Information hereunder, it shall be the burden of the receiving party to show both that such contested information or matter is not Proprietary Information within the meaning of this Agreement and that it does not constitute a trade secret under the Uniform Trade Secrets Act or successor or a similar law in effect in the State of (your state).',
    'iv. No delay or omission by either party in exercising any rights under this Agreement will operate as a waiver of that or any other right. A waiver or consent given by either party on any one occasion is effective only in that instance and will not be construed as a bar to or waiver of any right on any other occasion.',
    'v. This Agreement shall be binding upon and will inure to the benefit of the parties hereto and their respective successors and assigns.',
    'vi. This Agreement is governed by and will be construed in accordance with the laws of the State of New York, and the courts of New York shall be the exclusive forum.',
    'vii. This Agreement is in addition to any prior written agreement between Party 1 and Party 2 relating to the subject matter of this agreement; in the event of any disparity or conflict between the provision of such agreements, the provision which is more protective of Proprietary Information shall control. This Agreement may not be modified, in whole or in part, except by an agreement in writing signed by Party 1 and Party 2. IN WITNESS WHEREOF, the parties have executed this Agreement as of the date first above written. Omdena Inc. MUTUAL NON-DISCLOSURE AGREEMENT By: Rudradeb Mitra (Party 2 - Name & Signature here) Founder Omdena Inc.',
    '**Filled Data:**',
    '{',
    '"Annuities": 110.0,',
    '"Funeralwishes": "Thank you!",',
    '"SubExecAppointment": "fff",',
    '"Burial": "",',
    '"GiftOver": "",',
    '"GiftSpecificProperty": "",',
    '"SpecificBequests": "",',
    '"AgeofReceipt": "",',
    '"LifeInterests": "",',
    '"GuardianNeeded": "",',
    '"SpecificChattels": "",',
    '"TESTATORFULLNAME": "",',
    '"PeriodofGifting": "",',
    '"FamilyTrust": ""',
    '}',
    '**Cleaned-up Document:**',
    'Since no relevant placeholders from the document correlate with the filled data provided, the content remains the same as the original document. The placeholders for the date, Party2''s name, address, and registration number (${agreement-date}, ${Party2}, ${address}, ${register0number}) need details for completion. Further information is required to replace these placeholders accurately.'
)

# Overwrite the text of the first 37 paragraphs in place.
for ($i = 0; $i -lt $newTexts.Count; $i++) {
    $p = $d.Paragraphs.Item($i + 1)
    $p.Range.Text = $newTexts[$i]
}

# Remove the leftover trailing paragraphs (old content beyond the new 37).
$keepCount = $newTexts.Count
if ($d.Paragraphs.Count -gt $keepCount) {
    $startPara = $d.Paragraphs.Item($keepCount + 1)
    $endPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $trailingRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $trailingRange.Delete()
}

Write-Output $d.Paragraphs.Count
